$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "71.075.34"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.810.43"
$ws.Range("E3").Value = "  -0.89%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "700.96"
$ws.Range("E5").Value = "  -0.26%  "

# Row 6 - Solana
$ws.Range("D6").Value = "172.17"
$ws.Range("E6").Value = "  -0.39%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.809.05"
$ws.Range("E7").Value = "  -0.86%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +0.06%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.51%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "7.51"
$ws.Range("E11").Value = "  +2.08%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.38%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -1.11%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "36.02"
$ws.Range("E14").Value = "  -1.15%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.454.59"
$ws.Range("E15").Value = "  -0.92%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.882.03"
$ws.Range("E16").Value = "  +3.18%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "71.148.80"
$ws.Range("E17").Value = "  +0.24%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "17.50"
$ws.Range("E18").Value = "  +0.89%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "7.14"
$ws.Range("E19").Value = "  -0.39%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -0.55%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "512.64"
$ws.Range("E21").Value = "  +4.15%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "10.65"
$ws.Range("E22").Value = "  -0.29%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "0.715"
$ws.Range("E23").Value = "  -0.08%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "84.04"
$ws.Range("E24").Value = "  -1.19%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  -1.50%  "

# Row 26 - WrappedeETH
$ws.Range("D26").Value = "3.963.41"
$ws.Range("E26").Value = "  -0.98%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "12.08"
$ws.Range("E27").Value = "  -0.35%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "10.42"
$ws.Range("E28").Value = "  -1.04%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.09%  "

# Row 30 - Fetch.AI
$ws.Range("E30").Value = "  -3.42%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "3.03"
$ws.Range("E31").Value = "  -4.26%  "

# Row 32 - NEARProtocol
$ws.Range("D32").Value = "7.42"
$ws.Range("E32").Value = "  -1.06%  "

# Row 33 - ImmutableX
$ws.Range("D33").Value = "2.24"
$ws.Range("E33").Value = "  -1.53%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "29.07"
$ws.Range("E34").Value = "  -1.17%  "

# Row 35 - Kaspa
$ws.Range("D35").Value = "0.174"
$ws.Range("E35").Value = "  -4.64%  "

# Row 36 - Aptos
$ws.Range("E36").Value = "  +0.29%  "

# Row 37 - RenzoRestakedETH
$ws.Range("D37").Value = "3.773.38"
$ws.Range("E37").Value = "  -0.80%  "

# Row 38 - Binance-PegBSC-USD
$ws.Range("D38").Value = "0.997"
$ws.Range("E38").Value = "  -0.28%  "

# Row 39 - Hedera
$ws.Range("E39").Value = "  -1.98%  "

# Row 40 - Stacks
$ws.Range("D40").Value = "2.38"
$ws.Range("E40").Value = "  +0.65%  "

# Row 41 - Filecoin
$ws.Range("D41").Value = "6.01"
$ws.Range("E41").Value = "  -0.28%  "

# Row 42 - Mantle
$ws.Range("E42").Value = "  -0.68%  "

# Row 43 - dogwifhat
$ws.Range("E43").Value = "  -0.32%  "

# Row 45 - Monero
$ws.Range("D45").Value = "173.05"
$ws.Range("E45").Value = "  +5.88%  "

# Row 47 - FLOKI
$ws.Range("D47").Value = "0.000313"
$ws.Range("E47").Value = "  +1.42%  "

# Row 48 - was OKB, now Bittensor
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "431.59"
$ws.Range("E48").Value = "  +4.69%  "

# Row 49 - was Bittensor, now OKB
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "49.40"
$ws.Range("E49").Value = "  +1.44%  "

# Row 50 - Cosmos
$ws.Range("D50").Value = "8.63"
$ws.Range("E50").Value = "  +0.13%  "

# Row 51 - ONDO
$ws.Range("E51").Value = "  +0.72%  "
